$d = $word.ActiveDocument

$replacements = @(
    @("2025-03-14 Friday", "2025-03-15 Saturday"),
    @("641×8=5128", "859×9=7731"),
    @("657×9=5913", "283×4=1132"),
    @("772×2=1544", "359×7=2513"),
    @("997×2=1994", "604×8=4832"),
    @("353×2=706", "175×4=700"),
    @("129×8=1032", "430×2=860"),
    @("114×3=342", "268×9=2412"),
    @("507×3=1521", "819×7=5733"),
    @("538×4=2152", "193×3=579"),
    @("291×4=1164", "938×6=5628"),
    @("652×7=4564", "655×7=4585"),
    @("619×5=3095", "912×3=2736"),
    @("429×5=2145", "662×5=3310"),
    @("580×5=2900", "688×8=5504"),
    @("950×5=4750", "285×2=570"),
    @("393×9=3537", "791×4=3164"),
    @("932×5=4660", "930×7=6510"),
    @("930×8=7440", "644×5=3220"),
    @("463×3=1389", "475×3=1425"),
    @("196×2=392", "761×7=5327"),
    @("136×5=680", "767×9=6903"),
    @("534×9=4806", "299×5=1495"),
    @("643×8=5144", "683×5=3415"),
    @("457×6=2742", "680×4=2720"),
    @("326×2=652", "376×4=1504")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
